$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 1 questions (rows 2-6) -----------------------------------------
# Row 2: append new file names to the existing list of subjects/variables files
$ws.Range("D2").Value = " what-are-the-subjects-and-variables, identify-response-and-explanatory-vars, rectangular-data-interpretation, identify-subjects-and-variables"

# Row 3: append new file names + mark solution type as schoice
$ws.Range("D3").Value = "what-type-of-variable, identify-sample-space-and-variable-type-coin, identify-sample-space-and-variable-type-dice"
$ws.Range("E3").Value = "schoice"

# Row 4: new question file added
$ws.Range("D4").Value = "interpretation-coffee"

# Row 6: new question file + solution type added
$ws.Range("D6").Value = "sampling-bias-which-type, scientific-reasoning-coin, reliability-vs-validity"
$ws.Range("E6").Value = "schoice"

# --- Week 9 / ANOVA questions (rows 34-36) --------------------------------
$ws.Range("D34").Value = "anova-two-way-which-code-interaction, anova-two-way-interaction-means, anova-one-way-purpose, anova-one-way-sided-test,"
$ws.Range("E34").Value = "schoice"

$ws.Range("D35").Value = "anova-fraction-of-var-explained, anova-sided-test, anova-test-stat, anova-assumptions"

$ws.Range("D36").Value = "anova-two-way-conclusion, anova-two-way-significant-terms,  anova-one-way-interpret"
$ws.Range("E36").Value = "schoice"

# --- Row height adjustments (content reflow after edits) -----------------
$ws.Rows.Item(6).RowHeight = 48
$ws.Rows.Item(7).RowHeight = 48
$ws.Rows.Item(10).RowHeight = 17
$ws.Rows.Item(11).RowHeight = 34
$ws.Rows.Item(12).RowHeight = 32
$ws.Rows.Item(14).RowHeight = 34
$ws.Rows.Item(17).RowHeight = 34
$ws.Rows.Item(18).RowHeight = 34
$ws.Rows.Item(19).RowHeight = 34
$ws.Rows.Item(21).RowHeight = 34
$ws.Rows.Item(23).RowHeight = 17
$ws.Rows.Item(24).RowHeight = 34
$ws.Rows.Item(25).RowHeight = 34
$ws.Rows.Item(26).RowHeight = 34
$ws.Rows.Item(27).RowHeight = 34
$ws.Rows.Item(28).RowHeight = 34

# --- Column width adjustments (widened "Name of file" column, etc.) ------
$ws.Columns.Item(2).ColumnWidth = 5.2567
$ws.Columns.Item(3).ColumnWidth = 36.7567
$ws.Columns.Item(4).ColumnWidth = 107.09
$ws.Columns.Item(5).ColumnWidth = 8.7567

# --- Selection / view state ------------------------------------------------
$ws.Range("A5").Select()
$excel.ActiveWindow.Zoom = 100
